$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.40%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.60%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.702"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-11.60%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05970"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.675"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.46%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8713"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.87%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.75%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01065"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1,658.79%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1410"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03617"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.71%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07188"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.21%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03147"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.34%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001541"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.62%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005954"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.64%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.487"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.204"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.66%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.219"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.16%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3113"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.66%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-2.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.532"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.16%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04216"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.53%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.14%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.32%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-11.92%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-22.88%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03826"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.79%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'11.07%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.31%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002253"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01099"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.02%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005497"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.04%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.11%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1092"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'21.55%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002275"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5.58%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.11%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("E50").Style = "Normal"
